# Updates to mode init and HFO cap
# - Remove the modal-split initialization formulas for Road/Sea/Rail from column C
#   (replace with a flat "0.01" placeholder share) and move the real share
#   calculation out to column E under a new "ShareCalc" header.
# - Add reviewer comments (italic) explaining the rail-infeasibility issue and a
#   question about whether shares increased once Sweden was included.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # init_mode_mix
$ws2 = $wb.Worksheets.Item(2)   # init_fuel_mix

# --- init_mode_mix sheet -------------------------------------------------

# New header for the relocated share-calculation column.
$ws1.Range("E1").Value = "ShareCalc"
$ws1.Range("E1").Font.Bold = $true

# New review comments, added in shared-string order: G6 first, then G4.
$ws1.Range("G6").Value = "The shares in the model should be at least as big , as the ones denoted here"
$ws1.Range("G6").Font.Italic = $true

$ws1.Range("G4").Value = "maybe we have increased the shares by including sweden?"
$ws1.Range("G4").Font.Italic = $true

# Move the share formulas (previously in column C) over to column E.
$ws1.Range("E2").Formula = "=B2/SUM(`$B`$2:`$B`$4)"
$ws1.Range("E2").NumberFormat = "0.00"

$ws1.Range("E3").Formula = "=B3/SUM(`$B`$2:`$B`$4)"
$ws1.Range("E3").NumberFormat = "0.00"

$ws1.Range("E4").Formula = "=B4/SUM(`$B`$2:`$B`$4)"
$ws1.Range("E4").NumberFormat = "0.00"

# Column C becomes a flat placeholder share ("0.01") stored as literal text,
# not a number and not a formula. Enter it as a formula first so the text
# isn't re-parsed as a numeric literal, then paste-special as values only,
# then reset the cell style back to the sheet default.
foreach ($cellRef in @("C2", "C3", "C4")) {
    $cell = $ws1.Range($cellRef)
    $cell.Formula = "=""0.01"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $cell.Style = "Normal"
}

# Make init_mode_mix the active sheet/tab with C5 selected.
$ws1.Activate()
$ws1.Range("C5").Select()

# Match the new page setup (portrait, A4-ish paper size 9).
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
